$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '91.522.22'
$ws.Range("E2").Value = '  +0.98%  '
$ws.Range("D3").Value = '3.129.69'
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.59'
$ws.Range("E5").Value = '  -0.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '618.57'
$ws.Range("E6").Value = '  -0.87%  '
$ws.Range("E7").Value = '  -4.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.388'
$ws.Range("E8").Value = '  +4.56%  '
$ws.Range("D10").Value = '3.128.48'
$ws.Range("E10").Value = '  +0.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.749'
$ws.Range("E11").Value = '  -0.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.204'
$ws.Range("E12").Value = '  +0.41%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000254'
$ws.Range("E13").Value = '  +0.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.15'
$ws.Range("E14").Value = '  -0.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.61'
$ws.Range("E15").Value = '  +2.11%  '
$ws.Range("D16").Value = '91.337.65'
$ws.Range("E16").Value = '  +0.89%  '
$ws.Range("E17").Value = '  +0.87%  '
$ws.Range("D18").Value = '3.145.80'
$ws.Range("E18").Value = '  +1.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.76'
$ws.Range("E19").Value = '  -1.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.96'
$ws.Range("E20").Value = '  +4.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.93'
$ws.Range("E21").Value = '  +2.57%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '456.47'
$ws.Range("E22").Value = '  +1.92%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0000202'
$ws.Range("E23").Value = '  -4.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.22'
$ws.Range("E24").Value = '  +1.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.93'
$ws.Range("E25").Value = '  +0.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '88.83'
$ws.Range("E26").Value = '  -4.96%  '
$ws.Range("E27").Value = '  -0.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.153'
$ws.Range("E28").Value = '  +38.02%  '
$ws.Range("D29").Value = '3.312.58'
$ws.Range("E29").Value = '  +1.54%  '
$ws.Range("E30").Value = '  -0.14%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.232'
$ws.Range("E31").Value = '  +3.68%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.167'
$ws.Range("E32").Value = '  -5.50%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.37'
$ws.Range("E33").Value = '  +2.58%  '
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.175'
$ws.Range("E34").Value = '  +10.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '26.43'
$ws.Range("E35").Value = '  -0.70%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.44'
$ws.Range("E36").Value = '  -2.96%  '
$ws.Range("E37").Value = '  +2.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.97'
$ws.Range("E38").Value = '  -5.82%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '492.53'
$ws.Range("E39").Value = '  -0.33%  '
$ws.Range("E40").Value = '  +1.90%  '
$ws.Range("E41").Value = '  +5.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.39'
$ws.Range("E42").Value = '  -6.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.15'
$ws.Range("E43").Value = '  +0.17%  '
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.714'
$ws.Range("E45").Value = '  -28.46%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.709'
$ws.Range("E46").Value = '  +3.29%  '
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.93'
$ws.Range("E47").Value = '  +1.45%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '156.28'
$ws.Range("E48").Value = '  -0.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.36'
$ws.Range("E49").Value = '  +1.45%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.47'
$ws.Range("E50").Value = '  -2.24%  '
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0327'
$ws.Range("E51").Value = '  +5.43%  '
